$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "301.81"
Set-TextValue $ws.Range("E2") "-1.49%"
Set-TextValue $ws.Range("F2") "25-1-2023"
Set-TextValue $ws.Range("G2") "0"

# Row 3
Set-TextValue $ws.Range("D3") "35.54"
Set-TextValue $ws.Range("E3") "-2.04%"
Set-TextValue $ws.Range("F3") "25-1-2023"
Set-TextValue $ws.Range("G3") "0"

# Row 4
Set-TextValue $ws.Range("D4") "4.945"
Set-TextValue $ws.Range("E4") "-2.64%"
Set-TextValue $ws.Range("F4") "25-1-2023"
Set-TextValue $ws.Range("G4") "0"

# Row 5
Set-TextValue $ws.Range("D5") "0.07815"
Set-TextValue $ws.Range("E5") "-2.74%"
Set-TextValue $ws.Range("F5") "25-1-2023"
Set-TextValue $ws.Range("G5") "0"

# Row 6
Set-TextValue $ws.Range("D6") "1.867"
Set-TextValue $ws.Range("E6") "-14.18%"
Set-TextValue $ws.Range("F6") "25-1-2023"
Set-TextValue $ws.Range("G6") "0"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D7") "7.763"
Set-TextValue $ws.Range("E7") "-3.12%"
Set-TextValue $ws.Range("F7") "25-1-2023"
Set-TextValue $ws.Range("G7") "0"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D8") "2.932"
Set-TextValue $ws.Range("E8") "7.30%"
Set-TextValue $ws.Range("F8") "25-1-2023"
Set-TextValue $ws.Range("G8") "0"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9215"
Set-TextValue $ws.Range("E9") "-0.76%"
Set-TextValue $ws.Range("F9") "25-1-2023"
Set-TextValue $ws.Range("G9") "0"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1071"
Set-TextValue $ws.Range("E10") "8.73%"
Set-TextValue $ws.Range("F10") "25-1-2023"
Set-TextValue $ws.Range("G10") "0"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1837"
Set-TextValue $ws.Range("E11") "-2.39%"
Set-TextValue $ws.Range("F11") "25-1-2023"
Set-TextValue $ws.Range("G11") "0"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09275"
Set-TextValue $ws.Range("E12") "1.38%"
Set-TextValue $ws.Range("F12") "25-1-2023"
Set-TextValue $ws.Range("G12") "0"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03623"
Set-TextValue $ws.Range("E13") "-0.04%"
Set-TextValue $ws.Range("F13") "25-1-2023"
Set-TextValue $ws.Range("G13") "0"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09874"
Set-TextValue $ws.Range("E14") "-0.39%"
Set-TextValue $ws.Range("F14") "25-1-2023"
Set-TextValue $ws.Range("G14") "0"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001424"
Set-TextValue $ws.Range("E15") "-1.01%"
Set-TextValue $ws.Range("F15") "25-1-2023"
Set-TextValue $ws.Range("G15") "0"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005866"
Set-TextValue $ws.Range("E16") "3.12%"
Set-TextValue $ws.Range("F16") "25-1-2023"
Set-TextValue $ws.Range("G16") "0"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.455"
Set-TextValue $ws.Range("E17") "-0.11%"
Set-TextValue $ws.Range("F17") "25-1-2023"
Set-TextValue $ws.Range("G17") "0"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D18") "4.068"
Set-TextValue $ws.Range("E18") "-2.12%"
Set-TextValue $ws.Range("F18") "25-1-2023"
Set-TextValue $ws.Range("G18") "0"

# Row 19
Set-TextValue $ws.Range("D19") "0.3430"
Set-TextValue $ws.Range("F19") "25-1-2023"
Set-TextValue $ws.Range("G19") "0"

# Row 20
Set-TextValue $ws.Range("D20") "0.1309"
Set-TextValue $ws.Range("E20") "5.46%"
Set-TextValue $ws.Range("F20") "25-1-2023"
Set-TextValue $ws.Range("G20") "0"

# Row 21
Set-TextValue $ws.Range("D21") "5.106"
Set-TextValue $ws.Range("E21") "0.85%"
Set-TextValue $ws.Range("F21") "25-1-2023"
Set-TextValue $ws.Range("G21") "0"

# Row 22
Set-TextValue $ws.Range("D22") "0.2195"
Set-TextValue $ws.Range("E22") "-6.00%"
Set-TextValue $ws.Range("F22") "25-1-2023"
Set-TextValue $ws.Range("G22") "0"

# Row 23
Set-TextValue $ws.Range("D23") "0.04502"
Set-TextValue $ws.Range("E23") "-1.84%"
Set-TextValue $ws.Range("F23") "25-1-2023"
Set-TextValue $ws.Range("G23") "0"

# Row 24
Set-TextValue $ws.Range("D24") "0.001217"
Set-TextValue $ws.Range("E24") "-1.82%"
Set-TextValue $ws.Range("F24") "25-1-2023"
Set-TextValue $ws.Range("G24") "0"

# Row 25
Set-TextValue $ws.Range("D25") "0.004639"
Set-TextValue $ws.Range("E25") "-2.29%"
Set-TextValue $ws.Range("F25") "25-1-2023"
Set-TextValue $ws.Range("G25") "0"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001251"
Set-TextValue $ws.Range("E26") "-3.80%"
Set-TextValue $ws.Range("F26") "25-1-2023"
Set-TextValue $ws.Range("G26") "0"

# Row 27
Set-TextValue $ws.Range("D27") "0.0004460"
Set-TextValue $ws.Range("E27") "-0.88%"
Set-TextValue $ws.Range("F27") "25-1-2023"
Set-TextValue $ws.Range("G27") "0"

# Row 28
Set-TextValue $ws.Range("F28") "25-1-2023"
Set-TextValue $ws.Range("G28") "0"

# Row 29
Set-TextValue $ws.Range("F29") "25-1-2023"
Set-TextValue $ws.Range("G29") "0"

# Row 30
Set-TextValue $ws.Range("F30") "25-1-2023"
Set-TextValue $ws.Range("G30") "0"

# Row 31
Set-TextValue $ws.Range("F31") "25-1-2023"
Set-TextValue $ws.Range("G31") "0"

# Row 32
Set-TextValue $ws.Range("F32") "25-1-2023"
Set-TextValue $ws.Range("G32") "0"

# Row 33
Set-TextValue $ws.Range("F33") "25-1-2023"
Set-TextValue $ws.Range("G33") "0"

# Row 34
Set-TextValue $ws.Range("F34") "25-1-2023"
Set-TextValue $ws.Range("G34") "0"

# Row 35
Set-TextValue $ws.Range("F35") "25-1-2023"
Set-TextValue $ws.Range("G35") "0"

# Row 36
Set-TextValue $ws.Range("F36") "25-1-2023"
Set-TextValue $ws.Range("G36") "0"

# Row 37
Set-TextValue $ws.Range("F37") "25-1-2023"
Set-TextValue $ws.Range("G37") "0"

# Row 38
Set-TextValue $ws.Range("F38") "25-1-2023"
Set-TextValue $ws.Range("G38") "0"

# Row 39
Set-TextValue $ws.Range("D39") "0.01875"
Set-TextValue $ws.Range("E39") "-3.27%"
Set-TextValue $ws.Range("F39") "25-1-2023"
Set-TextValue $ws.Range("G39") "0"

# Row 40
Set-TextValue $ws.Range("D40") "0.04686"
Set-TextValue $ws.Range("E40") "-4.62%"
Set-TextValue $ws.Range("F40") "25-1-2023"
Set-TextValue $ws.Range("G40") "0"

# Row 41
Set-TextValue $ws.Range("D41") "0.007564"
Set-TextValue $ws.Range("E41") "-2.92%"
Set-TextValue $ws.Range("F41") "25-1-2023"
Set-TextValue $ws.Range("G41") "0"

# Row 42
Set-TextValue $ws.Range("D42") "0.009703"
Set-TextValue $ws.Range("E42") "24.15%"
Set-TextValue $ws.Range("F42") "25-1-2023"
Set-TextValue $ws.Range("G42") "0"

# Row 43
Set-TextValue $ws.Range("D43") "0.1330"
Set-TextValue $ws.Range("E43") "-4.77%"
Set-TextValue $ws.Range("F43") "25-1-2023"
Set-TextValue $ws.Range("G43") "0"

# Row 44
Set-TextValue $ws.Range("D44") "0.002116"
Set-TextValue $ws.Range("E44") "-0.68%"
Set-TextValue $ws.Range("F44") "25-1-2023"
Set-TextValue $ws.Range("G44") "0"

# Row 45
Set-TextValue $ws.Range("D45") "0.01115"
Set-TextValue $ws.Range("E45") "-1.48%"
Set-TextValue $ws.Range("F45") "25-1-2023"
Set-TextValue $ws.Range("G45") "0"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006148"
Set-TextValue $ws.Range("E46") "-1.70%"
Set-TextValue $ws.Range("F46") "25-1-2023"
Set-TextValue $ws.Range("G46") "0"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.03%"
Set-TextValue $ws.Range("F47") "25-1-2023"
Set-TextValue $ws.Range("G47") "0"

# Row 48
Set-TextValue $ws.Range("D48") "64.48"
Set-TextValue $ws.Range("E48") "138.33%"
Set-TextValue $ws.Range("F48") "25-1-2023"
Set-TextValue $ws.Range("G48") "0"

# Row 49
Set-TextValue $ws.Range("D49") "0.001301"
Set-TextValue $ws.Range("E49") "-27.70%"
Set-TextValue $ws.Range("F49") "25-1-2023"
Set-TextValue $ws.Range("G49") "0"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "0.03%"
Set-TextValue $ws.Range("F50") "25-1-2023"
Set-TextValue $ws.Range("G50") "0"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "0.03%"
Set-TextValue $ws.Range("F51") "25-1-2023"
Set-TextValue $ws.Range("G51") "0"
